# DebugDataTransfer.xlsx update: add new pipeline signals (0x18-0x1D) and
# rename the old "pipeline_jmp" entry (0x16) into the wider
# "pipeline_jmp_condl_rel_dests_cond_out" signal, per the register map.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24 (addr 0x16): rename signal + widen its bit width ------------
# Columns B/G already hold the address "0x16" and stay untouched.
$ws.Range("E24").Value = "pipeline_jmp_condl_rel_dests_cond_out"
$ws.Range("H24").Value = "7Bit"
$ws.Range("J24").Value = "pipeline_jmp_condl_rel_dests_cond_out"

# --- New rows 26-31: newly documented pipeline signals -------------------
$newRows = @(
    @{ Row = 26; Addr = "0x18"; Name = "pipeline_immediate_out";                     Width = "16 Bit" },
    @{ Row = 27; Addr = "0x19"; Name = "pipeline_write_address_out";                 Width = "4 Bit"  },
    @{ Row = 28; Addr = "0x1A"; Name = "pipeline_whb_wlb_out";                       Width = "2 Bit"  },
    @{ Row = 29; Addr = "0x1B"; Name = "pipeline_write_data_sel_out";                Width = "2 Bit"  },
    @{ Row = 30; Addr = "0x1C"; Name = "pipeline_RAM_src_read_write_bankid_out";     Width = "7 Bit"  },
    @{ Row = 31; Addr = "0x1D"; Name = "pipeline_is_alu_ram_gpu_op_out";             Width = "3 Bit"  }
)

foreach ($r in $newRows) {
    $ws.Range("B$($r.Row)").Value = $r.Addr
    $ws.Range("E$($r.Row)").Value = $r.Name
    $ws.Range("G$($r.Row)").Value = $r.Addr
    $ws.Range("H$($r.Row)").Value = $r.Width
    $ws.Range("J$($r.Row)").Value = $r.Name
}

# --- View state: move the selection like the author left it -------------
$ws.Range("D58").Select()
